$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 7).Value = 1.75
$ws.Cells.Item(3, 8).Value = 4.1
$ws.Cells.Item(3, 9).Value = 4.1
$ws.Cells.Item(3, 10).Value = 2.25
$ws.Cells.Item(3, 11).Value = 2.6
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(3, 13).Value = 1.02
$ws.Cells.Item(3, 14).Value = 21
$ws.Cells.Item(3, 15).Value = 1.11
$ws.Cells.Item(3, 16).Value = 6.5
$ws.Cells.Item(3, 17).Value = 1.4
$ws.Cells.Item(3, 18).Value = 2.88
$ws.Cells.Item(3, 19).Value = 1.22
$ws.Cells.Item(3, 20).Value = 4
$ws.Cells.Item(3, 21).Value = 1.4
$ws.Cells.Item(3, 22).Value = 2.75
$ws.Cells.Item(3, 23).Value = 13
$ws.Cells.Item(3, 24).Value = 12
$ws.Cells.Item(3, 25).Value = 9
$ws.Cells.Item(3, 26).Value = 17
$ws.Cells.Item(3, 27).Value = 12
$ws.Cells.Item(3, 28).Value = 17
$ws.Cells.Item(3, 29).Value = 21
$ws.Cells.Item(3, 30).Value = 8.5
$ws.Cells.Item(3, 31).Value = 11
$ws.Cells.Item(3, 32).Value = 29
$ws.Cells.Item(3, 33).Value = 81
$ws.Cells.Item(3, 34).Value = 21
$ws.Cells.Item(3, 35).Value = 26
$ws.Cells.Item(3, 36).Value = 15
$ws.Cells.Item(3, 37).Value = 41
$ws.Cells.Item(3, 38).Value = 26
$ws.Cells.Item(3, 39).Value = 26
$ws.Cells.Item(3, 40).Value = 4.33
$ws.Cells.Item(3, 41).Value = 9
$ws.Cells.Item(3, 42).Value = 15
$ws.Cells.Item(3, 43).Value = 23
$ws.Cells.Item(3, 44).Value = 41
$ws.Cells.Item(3, 45).Value = 67
$ws.Cells.Item(3, 46).Value = 4
$ws.Cells.Item(3, 47).Value = 7
$ws.Cells.Item(3, 48).Value = 41
$ws.Cells.Item(3, 49).Value = 6.5
$ws.Cells.Item(3, 50).Value = 19
$ws.Cells.Item(3, 51).Value = 21
$ws.Cells.Item(3, 52).Value = 51
$ws.Cells.Item(3, 53).Value = 51
$ws.Cells.Item(3, 54).Value = 101
$ws.Cells.Item(3, 55).Value = 251

# Row 4
$ws.Cells.Item(4, 7).Value = 3.7
$ws.Cells.Item(4, 8).Value = 3.6
$ws.Cells.Item(4, 9).Value = 1.95
$ws.Cells.Item(4, 10).Value = 3.75
$ws.Cells.Item(4, 11).Value = 2.3
$ws.Cells.Item(4, 12).Value = 2.6
$ws.Cells.Item(4, 13).Value = 1.03
$ws.Cells.Item(4, 14).Value = 15
$ws.Cells.Item(4, 15).Value = 1.18
$ws.Cells.Item(4, 16).Value = 4.5
$ws.Cells.Item(4, 17).Value = 1.65
$ws.Cells.Item(4, 18).Value = 2.2
$ws.Cells.Item(4, 19).Value = 1.3
$ws.Cells.Item(4, 20).Value = 3.4
$ws.Cells.Item(4, 21).Value = 1.57
$ws.Cells.Item(4, 22).Value = 2.25
$ws.Cells.Item(4, 23).Value = 15
$ws.Cells.Item(4, 24).Value = 21
$ws.Cells.Item(4, 25).Value = 13
$ws.Cells.Item(4, 26).Value = 41
$ws.Cells.Item(4, 27).Value = 26
$ws.Cells.Item(4, 28).Value = 29
$ws.Cells.Item(4, 29).Value = 15
$ws.Cells.Item(4, 30).Value = 7
$ws.Cells.Item(4, 31).Value = 12
$ws.Cells.Item(4, 32).Value = 41
$ws.Cells.Item(4, 33).Value = 126
$ws.Cells.Item(4, 34).Value = 9.5
$ws.Cells.Item(4, 35).Value = 11
$ws.Cells.Item(4, 36).Value = 9
$ws.Cells.Item(4, 37).Value = 19
$ws.Cells.Item(4, 38).Value = 15
$ws.Cells.Item(4, 39).Value = 21
$ws.Cells.Item(4, 40).Value = 5.5
$ws.Cells.Item(4, 41).Value = 19
$ws.Cells.Item(4, 42).Value = 23
$ws.Cells.Item(4, 43).Value = 51
$ws.Cells.Item(4, 44).Value = 67
$ws.Cells.Item(4, 45).Value = 126
$ws.Cells.Item(4, 46).Value = 3.4
$ws.Cells.Item(4, 47).Value = 7.5
$ws.Cells.Item(4, 48).Value = 41
$ws.Cells.Item(4, 49).Value = 4.33
$ws.Cells.Item(4, 50).Value = 10
$ws.Cells.Item(4, 51).Value = 19
$ws.Cells.Item(4, 52).Value = 34
$ws.Cells.Item(4, 53).Value = 51
$ws.Cells.Item(4, 54).Value = 101

# Row 5
$ws.Cells.Item(5, 7).Value = 2.25
$ws.Cells.Item(5, 8).Value = 3.25
$ws.Cells.Item(5, 9).Value = 3.2
$ws.Cells.Item(5, 10).Value = 2.88
$ws.Cells.Item(5, 11).Value = 2.3
$ws.Cells.Item(5, 12).Value = 3.4
$ws.Cells.Item(5, 13).Value = 1.03
$ws.Cells.Item(5, 14).Value = 15
$ws.Cells.Item(5, 15).Value = 1.18
$ws.Cells.Item(5, 16).Value = 4.5
$ws.Cells.Item(5, 17).Value = 1.62
$ws.Cells.Item(5, 18).Value = 2.25
$ws.Cells.Item(5, 19).Value = 1.3
$ws.Cells.Item(5, 20).Value = 3.4
$ws.Cells.Item(5, 21).Value = 1.5
$ws.Cells.Item(5, 22).Value = 2.5
$ws.Cells.Item(5, 23).Value = 11
$ws.Cells.Item(5, 24).Value = 13
$ws.Cells.Item(5, 25).Value = 9.5
$ws.Cells.Item(5, 26).Value = 21
$ws.Cells.Item(5, 27).Value = 17
$ws.Cells.Item(5, 28).Value = 21
$ws.Cells.Item(5, 29).Value = 13
$ws.Cells.Item(5, 30).Value = 6.5
$ws.Cells.Item(5, 31).Value = 11
$ws.Cells.Item(5, 32).Value = 34
$ws.Cells.Item(5, 33).Value = 101
$ws.Cells.Item(5, 34).Value = 13
$ws.Cells.Item(5, 35).Value = 19
$ws.Cells.Item(5, 36).Value = 12
$ws.Cells.Item(5, 37).Value = 34
$ws.Cells.Item(5, 38).Value = 23
$ws.Cells.Item(5, 39).Value = 26
$ws.Cells.Item(5, 40).Value = 4.5
$ws.Cells.Item(5, 41).Value = 12
$ws.Cells.Item(5, 42).Value = 19
$ws.Cells.Item(5, 43).Value = 41
$ws.Cells.Item(5, 44).Value = 51
$ws.Cells.Item(5, 45).Value = 101
$ws.Cells.Item(5, 46).Value = 3.4
$ws.Cells.Item(5, 47).Value = 7
$ws.Cells.Item(5, 48).Value = 41
$ws.Cells.Item(5, 49).Value = 5
$ws.Cells.Item(5, 50).Value = 15
$ws.Cells.Item(5, 51).Value = 21
$ws.Cells.Item(5, 52).Value = 51
$ws.Cells.Item(5, 53).Value = 51
$ws.Cells.Item(5, 54).Value = 126

